# Apply the "Append: 2026-01-09 01:27 JST" update to the "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-09 01:27:42"

# Update the "取得日時" (fetched at) timestamp for every data row (2..13).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Row 3 content changes: title, priority score and skill summary were
# revised (the "購入bot" / "★bot " wording was dropped, and the score
# dropped from 235 to 120).
$ws.Range("B3").Value = "初回 急募 自動カートインツール 開発のプロフェッショナルを探しています"
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = "◆ツール,開発"

# Column B got a bit narrower (43 -> 38 characters).
$ws.Columns.Item(2).ColumnWidth = 37.17
